$d = $word.ActiveDocument

# --- 1. Remove the _GoBack bookmark after "...) a plot." ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Append new paragraphs at the end of the document (after the
#        existing trailing empty paragraph, before the section break) ---
function Insert-ParaXml([string]$innerXml) {
    $doc = $word.ActiveDocument
    $lastPara = $doc.Paragraphs.Last
    $r = $lastPara.Range
    $pkg = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>" + "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" + "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" + "<pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body>" + $innerXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $r.InsertXML($pkg)
}

# First, create a fresh paragraph after the current last (empty) paragraph;
# each Insert-ParaXml call below then lands its content directly in front of
# that still-empty trailing paragraph, in document order.
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve">How to publish </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>Quarto using terminal:</w:t>
      </w:r>
    </w:p>'
Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>quarto publish quarto-pub</w:t>
      </w:r>
    </w:p>'
Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
    </w:p>'
Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>library(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>gt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
    </w:p>'
Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>exibble</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve"> %&gt;%</w:t>
      </w:r>
    </w:p>'
Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>gt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>) %&gt;%</w:t>
      </w:r>
    </w:p>'
Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>fmt_</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>number</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>'
Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve">    columns = num,</w:t>
      </w:r>
    </w:p>'
Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve">    decimals = 3,</w:t>
      </w:r>
    </w:p>'
Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>use_seps</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve"> = FLASE</w:t>
      </w:r>
    </w:p>'
Insert-ParaXml '<w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="NanumGothic"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve">  )</w:t>
      </w:r>
      <w:bookmarkEnd w:id="0"/>
    </w:p>'

# The InsertParagraphAfter() call above leaves one extra, now-redundant
# empty paragraph at the very end (each Insert-ParaXml shifted its own
# content in front of it). Remove that leftover paragraph mark so the
# last inserted paragraph ("  )") is immediately followed by the section
# break, matching the target document.
$finalPara = $d.Paragraphs.Last
$priorPara = $finalPara.Previous()
$cleanupRange = $d.Range($priorPara.Range.End - 1, $finalPara.Range.End)
$cleanupRange.Delete()

Write-Output "done"
